# Applies the "A little bit of data updated" commit to Data.xlsx.
#
# Sheet "COMPOUNDS" (sheet1): species list changed from
#   Water / Nitrogen / Oxygen / Cooling Water  ->  Water / Nitrogen / Methane / Hydrogen
# and the chem-formula row switched from LaTeX ($H_2O$ ...) to plain text
# (H20, N2, CH4, H2) plus a new 4th compound column (HCN).
#
# Sheet "UNIT OPERATIONS" (sheet2): unit list changed from
#   Reactor / Column / Sparger / Crystalizer -> Reactor / NH3Absorber / HCNAbsorber / Distillation
# several numeric columns were cleared out, "Inlet/Outlet flow" became
# "Height"/"Radius", "RedBull" became "En", "reactir simulation time" became
# "unit simulation time", and two new rows (Lebensdauer/lspan with a formula,
# and VolumeVoid/epsilon) were appended.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: COMPOUNDS
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("COMPOUNDS")

# Row 3: compound names - Oxygen -> Methane, Cooling Water -> Hydrogen
$ws1.Range("E3").Value = 'Methane'
$ws1.Range("F3").Value = 'Hydrogen'

# Row 4: label + chemical formulas (switch from LaTeX to plain identifiers,
# add a 4th/5th compound formula)
$ws1.Range("A4").Value = 'ChemFormula'
$ws1.Range("C4").Value = 'H20'
$ws1.Range("D4").Value = 'N2'
$ws1.Range("E4").Value = 'CH4'
$ws1.Range("F4").Value = 'H2'
$ws1.Range("G4").Value = 'HCN'

$ws1.Range("H4").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: UNIT OPERATIONS
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("UNIT OPERATIONS")

# Row 3: unit names
$ws2.Range("D3").Value = 'NH3Absorber'
$ws2.Range("E3").Value = 'HCNAbsorber'
$ws2.Range("F3").Value = 'Distillation'

# Row 4 (Volume/V): clear out the old per-unit numbers
$ws2.Range("C4:E4").ClearContents()

# Row 5: "Inlet flow"/"Qin" -> "Height"/"h", clear old numbers
$ws2.Range("A5").Value = 'Height'
$ws2.Range("B5").Value = 'h'
$ws2.Range("C5:F5").ClearContents()

# Row 6: "Outlet Flow"/"Qout" -> "Radius"/"rad", clear old numbers
$ws2.Range("A6").Value = 'Radius'
$ws2.Range("B6").Value = 'rad'
$ws2.Range("C6:F6").ClearContents()

# Row 7 (Optimal conversion/Xopt): clear out the old per-unit numbers
$ws2.Range("C7:F7").ClearContents()

# Row 8: "RedBull" -> "En", clear the stored energy value (style stays)
$ws2.Range("B8").Value = 'En'
$ws2.Range("E8").ClearContents()

# Row 9 (CAPEX/capex): clear stored value(s)
$ws2.Range("C9").ClearContents()
$ws2.Range("D9:F9").ClearContents()

# Row 10 (OPEX/opex): clear stored value
$ws2.Range("E10").ClearContents()

# Row 11 (Total Cost/totex): clear stored values
$ws2.Range("C11:D11").ClearContents()
$ws2.Range("G11").ClearContents()

# Row 12: "reactir simulation time" -> "unit simulation time"
$ws2.Range("A12").Value = 'unit simulation time'

# Row 13 (new): Lebensdauer / lspan, with a lifetime formula in G
$ws2.Range("A13").Value = 'Lebensdauer'
$ws2.Range("B13").Value = 'lspan'
$ws2.Range("G13").Formula = '=10*365*24*3600'

# Row 14 (new): VolumeVoid / epsilon
$ws2.Range("A14").Value = 'VolumeVoid'
$ws2.Range("B14").Value = 'epsilon'
$ws2.Range("F14").Value = 0.74

$ws2.Activate() | Out-Null
$ws2.Range("F15").Select() | Out-Null
